$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A196").Value = 195
$ws.Range("B196").Value = 1
$ws.Range("C196").Value = "2024-06-19 03:14:20"
$ws.Range("D196").Value = 200
$ws.Range("E196").Value = 20

$ws.Range("A197").Value = 196
$ws.Range("B197").Value = 2
$ws.Range("C197").Value = "2024-06-19 03:14:21"
$ws.Range("D197").Value = 200
$ws.Range("E197").Value = 3
